$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old content in A1
$ws.Range("A1").ClearContents()

# Write the new values - same shared string reused in A2 and A3
$ws.Range("A2").Value = "second test"
$ws.Range("A3").Value = "second test"

# AutoFit column A width (bestFit)
$ws.Columns.Item(1).AutoFit() | Out-Null

# Update the active selection
$ws.Range("G5").Select() | Out-Null
